$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$ws.Range('G2').Value = 'system, System, backup@backdoor.com'
$ws.Range('G3').Value = 'System, dnasr281@gmail.com'
$ws.Range('G4').Value = 'System, backup@backdoor.com'
$ws.Range('G5').Value = 'System, backup@backdoor.com'
$ws.Range('G6').Value = 'System, dnasr281@gmail.com'
$ws.Range('G7').Value = 'System, admin@admin.com'
$ws.Range('G10').Value = 'System, dnasr281@gmail.com'
$ws.Range('G12').Value = 'System, dnasr281@gmail.com'
$ws.Range('G13').Value = 'System, dnasr281@gmail.com'
$ws.Range('G14').Value = 'System, dnasr281@gmail.com'
$ws.Range('G15').Value = 'System, dnasr281@gmail.com'
$ws.Range('G29').Value = 'system, System, backup@backdoor.com'
$ws.Range('G30').Value = 'System, dnasr281@gmail.com'
$ws.Range('G32').Value = 'System, backup@backdoor.com'
$ws.Range('G33').Value = 'System, dnasr281@gmail.com'
$ws.Range('G34').Value = 'System, admin@admin.com'
$ws.Range('G37').Value = 'System, dnasr281@gmail.com'
$ws.Range('G39').Value = 'System, dnasr281@gmail.com'
$ws.Range('G40').Value = 'System, dnasr281@gmail.com'
$ws.Range('G41').Value = 'System, dnasr281@gmail.com'
$ws.Range('G42').Value = 'System, dnasr281@gmail.com'
$ws.Range('G56').Value = 'system, System, backup@backdoor.com'
$ws.Range('G57').Value = 'System, dnasr281@gmail.com'
$ws.Range('G58').Value = 'System, backup@backdoor.com'
$ws.Range('G59').Value = 'System, backup@backdoor.com'
$ws.Range('G60').Value = 'System, dnasr281@gmail.com'
$ws.Range('G61').Value = 'System, admin@admin.com'
$ws.Range('G64').Value = 'System, dnasr281@gmail.com'
$ws.Range('G66').Value = 'System, dnasr281@gmail.com'
$ws.Range('G67').Value = 'System, dnasr281@gmail.com'
$ws.Range('G68').Value = 'System, dnasr281@gmail.com'
$ws.Range('G69').Value = 'System, dnasr281@gmail.com'
$ws.Range('G84').Value = 'System, backup@backdoor.com'
$ws.Range('G85').Value = 'System, backup@backdoor.com'
$ws.Range('G86').Value = 'System, dnasr281@gmail.com'
$ws.Range('G87').Value = 'System, dnasr281@gmail.com'
$ws.Range('G88').Value = 'System, dnasr281@gmail.com'
$ws.Range('G89').Value = 'System, dnasr281@gmail.com'
$ws.Range('G90').Value = 'dnasr281@gmail.com, admin@admin.com'
$ws.Range('G93').Value = 'System, dnasr281@gmail.com'
$ws.Range('G95').Value = 'System, dnasr281@gmail.com'
$ws.Range('G110').Value = 'System, backup@backdoor.com'
$ws.Range('G111').Value = 'System, backup@backdoor.com'
$ws.Range('G112').Value = 'System, dnasr281@gmail.com'
$ws.Range('G113').Value = 'System, dnasr281@gmail.com'
$ws.Range('G114').Value = 'System, dnasr281@gmail.com'
$ws.Range('G115').Value = 'System, dnasr281@gmail.com'
$ws.Range('G116').Value = 'dnasr281@gmail.com, admin@admin.com'
$ws.Range('G119').Value = 'System, dnasr281@gmail.com'
$ws.Range('G121').Value = 'System, dnasr281@gmail.com'
$ws.Range('G136').Value = 'System, backup@backdoor.com'
$ws.Range('G137').Value = 'System, backup@backdoor.com'
$ws.Range('G138').Value = 'System, dnasr281@gmail.com'
$ws.Range('G139').Value = 'System, dnasr281@gmail.com'
$ws.Range('G140').Value = 'System, dnasr281@gmail.com'
$ws.Range('G141').Value = 'System, dnasr281@gmail.com'
$ws.Range('G142').Value = 'dnasr281@gmail.com, admin@admin.com'
$ws.Range('G145').Value = 'System, dnasr281@gmail.com'
$ws.Range('G147').Value = 'System, dnasr281@gmail.com'
